$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 3947.162
$ws.Cells.Item(69, 9).Value = 3400
$ws.Cells.Item(69, 10).Value = 8461.25
$ws.Cells.Item(69, 11).Value = 10200
$ws.Cells.Item(69, 12).Value = 25383.75
$ws.Cells.Item(69, 13).Value = -9326
$ws.Cells.Item(69, 14).Value = -27131.75
$ws.Cells.Item(72, 8).Value = 3947.162
$ws.Cells.Item(72, 9).Value = 3400
$ws.Cells.Item(72, 10).Value = 8461.25
$ws.Cells.Item(72, 11).Value = 30600
$ws.Cells.Item(72, 12).Value = 76151.25
$ws.Cells.Item(72, 13).Value = -26232
$ws.Cells.Item(72, 14).Value = -84887.25
$ws.Cells.Item(86, 8).Value = 102406.9
$ws.Cells.Item(86, 9).Value = 201918.2
$ws.Cells.Item(86, 10).Value = 2895.6
$ws.Cells.Item(86, 11).Value = 201918.2
$ws.Cells.Item(86, 12).Value = 2895.6
$ws.Cells.Item(86, 13).Value = -200795.2
$ws.Cells.Item(86, 14).Value = -5141.6
$ws.Cells.Item(89, 8).Value = 102406.9
$ws.Cells.Item(89, 9).Value = 201918.2
$ws.Cells.Item(89, 10).Value = 2895.6
$ws.Cells.Item(89, 11).Value = 1009591
$ws.Cells.Item(89, 12).Value = 14478
$ws.Cells.Item(89, 13).Value = -1003975
$ws.Cells.Item(89, 14).Value = -25710
$ws.Cells.Item(132, 8).Value = 2977761.8
$ws.Cells.Item(132, 9).Value = 4083317.2
$ws.Cells.Item(132, 10).Value = 1265.9231
$ws.Cells.Item(132, 11).Value = 12249951.6
$ws.Cells.Item(132, 12).Value = 3797.7693
$ws.Cells.Item(132, 13).Value = -12247421.6
$ws.Cells.Item(132, 14).Value = -8857.7693
$ws.Cells.Item(137, 8).Value = 1126.3125
$ws.Cells.Item(137, 9).Value = 793.72
$ws.Cells.Item(137, 10).Value = 2314.1428
$ws.Cells.Item(137, 11).Value = 2381.16
$ws.Cells.Item(137, 12).Value = 6942.428400000001
$ws.Cells.Item(137, 13).Value = 168.8400000000001
$ws.Cells.Item(137, 14).Value = -12042.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 671551.0600000001
$ws.Cells.Item(88, 9).Value = 1432849.2
$ws.Cells.Item(88, 10).Value = 5415.125
$ws.Cells.Item(88, 11).Value = 1432849.2
$ws.Cells.Item(88, 12).Value = 5415.125
$ws.Cells.Item(88, 13).Value = -1432443.2
$ws.Cells.Item(88, 14).Value = -6227.125
$ws.Cells.Item(91, 8).Value = 671551.0600000001
$ws.Cells.Item(91, 9).Value = 1432849.2
$ws.Cells.Item(91, 10).Value = 5415.125
$ws.Cells.Item(91, 11).Value = 1432849.2
$ws.Cells.Item(91, 12).Value = 5415.125
$ws.Cells.Item(91, 13).Value = -1431445.2
$ws.Cells.Item(91, 14).Value = -8223.125
$ws.Cells.Item(122, 8).Value = 1134.1154
$ws.Cells.Item(122, 9).Value = 860.44446
$ws.Cells.Item(122, 10).Value = 1749.875
$ws.Cells.Item(122, 11).Value = 2581.33338
$ws.Cells.Item(122, 12).Value = 5249.625
$ws.Cells.Item(122, 13).Value = -131.33338
$ws.Cells.Item(122, 14).Value = -10149.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2835.5
$ws.Cells.Item(86, 9).Value = 1900
$ws.Cells.Item(86, 10).Value = 3396.8
$ws.Cells.Item(86, 11).Value = 1900
$ws.Cells.Item(86, 12).Value = 3396.8
$ws.Cells.Item(86, 13).Value = -777
$ws.Cells.Item(86, 14).Value = -5642.8
$ws.Cells.Item(89, 8).Value = 2835.5
$ws.Cells.Item(89, 9).Value = 1900
$ws.Cells.Item(89, 10).Value = 3396.8
$ws.Cells.Item(89, 11).Value = 9500
$ws.Cells.Item(89, 12).Value = 16984
$ws.Cells.Item(89, 13).Value = -3884
$ws.Cells.Item(89, 14).Value = -28216
$ws.Cells.Item(134, 8).Value = 61171.117
$ws.Cells.Item(134, 9).Value = 112789.89
$ws.Cells.Item(134, 10).Value = 3100
$ws.Cells.Item(134, 11).Value = 338369.67
$ws.Cells.Item(134, 12).Value = 9300
$ws.Cells.Item(134, 13).Value = -335834.67
$ws.Cells.Item(134, 14).Value = -14370

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3626025.8
$ws.Cells.Item(31, 9).Value = 2785.7144
$ws.Cells.Item(31, 10).Value = 6669547.5
$ws.Cells.Item(31, 11).Value = 2785.7144
$ws.Cells.Item(31, 12).Value = 6669547.5
$ws.Cells.Item(31, 13).Value = -2490.7144
$ws.Cells.Item(31, 14).Value = -6670137.5
$ws.Cells.Item(34, 8).Value = 3626025.8
$ws.Cells.Item(34, 9).Value = 2785.7144
$ws.Cells.Item(34, 10).Value = 6669547.5
$ws.Cells.Item(34, 11).Value = 2785.7144
$ws.Cells.Item(34, 12).Value = 6669547.5
$ws.Cells.Item(34, 13).Value = -2583.7144
$ws.Cells.Item(34, 14).Value = -6669951.5
$ws.Cells.Item(134, 8).Value = 1150
$ws.Cells.Item(134, 9).Value = 955.55554
$ws.Cells.Item(134, 10).Value = 1733.3334
$ws.Cells.Item(134, 11).Value = 2866.66662
$ws.Cells.Item(134, 12).Value = 5200.0002
$ws.Cells.Item(134, 13).Value = -331.66662
$ws.Cells.Item(134, 14).Value = -10270.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1158.6842
$ws.Cells.Item(5, 9).Value = 219.66667
$ws.Cells.Item(5, 10).Value = 2003.8
$ws.Cells.Item(5, 11).Value = 659.00001
$ws.Cells.Item(5, 12).Value = 6011.4
$ws.Cells.Item(5, 13).Value = -547.00001
$ws.Cells.Item(5, 14).Value = -6235.4
$ws.Cells.Item(97, 8).Value = 398
$ws.Cells.Item(97, 9).Value = 397
$ws.Cells.Item(97, 10).Value = 399.2
$ws.Cells.Item(97, 11).Value = 1191
$ws.Cells.Item(97, 12).Value = 1197.6
$ws.Cells.Item(97, 13).Value = -695
$ws.Cells.Item(97, 14).Value = -2189.6
$ws.Cells.Item(132, 8).Value = 2098.625
$ws.Cells.Item(132, 9).Value = 995.875
$ws.Cells.Item(132, 10).Value = 2650
$ws.Cells.Item(132, 11).Value = 8962.875
$ws.Cells.Item(132, 12).Value = 23850
$ws.Cells.Item(132, 13).Value = -6432.875
$ws.Cells.Item(132, 14).Value = -28910
$ws.Cells.Item(135, 8).Value = 1158.6842
$ws.Cells.Item(135, 9).Value = 219.66667
$ws.Cells.Item(135, 10).Value = 2003.8
$ws.Cells.Item(135, 11).Value = 1977.00003
$ws.Cells.Item(135, 12).Value = 18034.2
$ws.Cells.Item(135, 13).Value = 557.9999699999998
$ws.Cells.Item(135, 14).Value = -23104.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 563
$ws.Cells.Item(16, 9).Value = 345
$ws.Cells.Item(16, 10).Value = 999
$ws.Cells.Item(16, 11).Value = 345
$ws.Cells.Item(16, 12).Value = 999
$ws.Cells.Item(16, 13).Value = -175
$ws.Cells.Item(16, 14).Value = -1339
$ws.Cells.Item(61, 8).Value = 15153191
$ws.Cells.Item(61, 9).Value = 1787.9445
$ws.Cells.Item(61, 10).Value = 83334504
$ws.Cells.Item(61, 11).Value = 1787.9445
$ws.Cells.Item(61, 12).Value = 83334504
$ws.Cells.Item(61, 13).Value = -1585.9445
$ws.Cells.Item(61, 14).Value = -83334908
$ws.Cells.Item(113, 8).Value = 15153191
$ws.Cells.Item(113, 9).Value = 1787.9445
$ws.Cells.Item(113, 10).Value = 83334504
$ws.Cells.Item(113, 11).Value = 1787.9445
$ws.Cells.Item(113, 12).Value = 83334504
$ws.Cells.Item(113, 13).Value = 382.0554999999999
$ws.Cells.Item(113, 14).Value = -83338844
$ws.Cells.Item(132, 8).Value = 34200.332
$ws.Cells.Item(132, 9).Value = 63734.668
$ws.Cells.Item(132, 10).Value = 4666
$ws.Cells.Item(132, 11).Value = 191204.004
$ws.Cells.Item(132, 12).Value = 13998
$ws.Cells.Item(132, 13).Value = -188674.004
$ws.Cells.Item(132, 14).Value = -19058
$ws.Cells.Item(136, 8).Value = 6167.75
$ws.Cells.Item(136, 9).Value = 11190.4
$ws.Cells.Item(136, 10).Value = 2580.1428
$ws.Cells.Item(136, 11).Value = 33571.2
$ws.Cells.Item(136, 12).Value = 7740.428400000001
$ws.Cells.Item(136, 13).Value = -31021.2
$ws.Cells.Item(136, 14).Value = -12840.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4371.4287
$ws.Cells.Item(62, 9).Value = 4275
$ws.Cells.Item(62, 10).Value = 4500
$ws.Cells.Item(62, 11).Value = 4275
$ws.Cells.Item(62, 12).Value = 4500
$ws.Cells.Item(62, 13).Value = -3651
$ws.Cells.Item(62, 14).Value = -5748
$ws.Cells.Item(65, 8).Value = 4371.4287
$ws.Cells.Item(65, 9).Value = 4275
$ws.Cells.Item(65, 10).Value = 4500
$ws.Cells.Item(65, 11).Value = 21375
$ws.Cells.Item(65, 12).Value = 22500
$ws.Cells.Item(65, 13).Value = -18255
$ws.Cells.Item(65, 14).Value = -28740
$ws.Cells.Item(81, 8).Value = 2704.158
$ws.Cells.Item(81, 9).Value = 2125.3635
$ws.Cells.Item(81, 10).Value = 3500
$ws.Cells.Item(81, 11).Value = 4250.727
$ws.Cells.Item(81, 12).Value = 7000
$ws.Cells.Item(81, 13).Value = -3189.727
$ws.Cells.Item(81, 14).Value = -9122
$ws.Cells.Item(84, 8).Value = 2704.158
$ws.Cells.Item(84, 9).Value = 2125.3635
$ws.Cells.Item(84, 10).Value = 3500
$ws.Cells.Item(84, 11).Value = 21253.635
$ws.Cells.Item(84, 12).Value = 35000
$ws.Cells.Item(84, 13).Value = -15949.635
$ws.Cells.Item(84, 14).Value = -45608
$ws.Cells.Item(110, 8).Value = 38780
$ws.Cells.Item(110, 9).Value = 0
$ws.Cells.Item(110, 10).Value = 38780
$ws.Cells.Item(110, 11).Value = 0
$ws.Cells.Item(110, 12).Value = 38780
$ws.Cells.Item(110, 14).Value = -46960
$ws.Cells.Item(126, 8).Value = 47625264
$ws.Cells.Item(126, 9).Value = 100011870
$ws.Cells.Item(126, 10).Value = 1076.5454
$ws.Cells.Item(126, 11).Value = 300035610
$ws.Cells.Item(126, 12).Value = 3229.6362
$ws.Cells.Item(126, 13).Value = -300033140
$ws.Cells.Item(126, 14).Value = -8169.6362
$ws.Cells.Item(132, 8).Value = 2160.5454
$ws.Cells.Item(132, 9).Value = 1695.6
$ws.Cells.Item(132, 10).Value = 3156.8572
$ws.Cells.Item(132, 11).Value = 5086.799999999999
$ws.Cells.Item(132, 12).Value = 9470.571599999999
$ws.Cells.Item(132, 13).Value = -2556.799999999999
$ws.Cells.Item(132, 14).Value = -14530.5716
$ws.Cells.Item(136, 8).Value = 3756.878
$ws.Cells.Item(136, 9).Value = 4231.5293
$ws.Cells.Item(136, 10).Value = 1451.4286
$ws.Cells.Item(136, 11).Value = 12694.5879
$ws.Cells.Item(136, 12).Value = 4354.2858
$ws.Cells.Item(136, 13).Value = -10144.5879
$ws.Cells.Item(136, 14).Value = -9454.2858
